$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column (C) for every data row (2..416)
# was bumped by one day: 45189 (2023-09-20) -> 45190 (2023-09-21)
$ws.Range("C2:C416").Value = 45190
